$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Duplicate the formatting of the last existing row (33) onto the new row (34)
$ws.Range("A33:E33").Copy()
$ws.Range("A34:E34").PasteSpecial(-4122)

# Fill in the new testscript's data (WAT46)
$ws.Range("A34").Value = "WAT46"
$ws.Range("C34").Value = "Verify that user should be able to view profile details by clicking on the profile pic button top right corner."
$ws.Range("B34").Value = "WAT-161"
$ws.Range("D34").Value = "Y"

# Match the workbook's final selection state
$ws.Range("E34").Select()
